$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- New data rows (5-7): Maven build-plugin tips ----

# Row 5: jdk1.8 compatibility
$ws.Range("A5").Value = 'Maven'
$ws.Range("B5").Value = 'jdk1.8 not compate'
$ws.Range("C5").Value = 'Maven assumes to build witih jdk1.5, add below code to pom for other jdk version:
  <build>
    <plugins>
   <plugin>
        <groupId>org.apache.maven.plugins</groupId>
        <artifactId>maven-compiler-plugin</artifactId>
        <version>3.3</version>
        <configuration>
            <source>1.8</source>
            <target>1.8</target>
        </configuration>
   </plugin>
    </plugins>
  </build>'

# Row 6: building a jar together with its dependencies
$ws.Range("A6").Value = 'Maven'
$ws.Range("B6").Value = 'Build with dependency'
$ws.Range("C6").Value = 'By default, Maven doesn''t bundle dependencies in the JAR file it builds, and you''re not providing them on the classpath when you''re trying to execute your JAR file at the command-line. This is why the Java VM can''t find the library class files when trying to execute your code.
You could manually specify the libraries on the classpath with the -cp parameter, or There is a Maven plugin called the maven-shade-plugin to do this. You need to register it in your POM, and it will automatically build an "uber-JAR" containing your classes and the classes for your library code too when you run mvn package:
  <build>
    <plugins>
      <plugin>
        <groupId>org.apache.maven.plugins</groupId>
        <artifactId>maven-shade-plugin</artifactId>
        <version>1.6</version>
        <executions>
          <execution>
            <phase>package</phase>
            <goals>
              <goal>shade</goal>
            </goals>
          </execution>
        </executions>
      </plugin>
    </plugins>
  </build>'

# Row 7: building a plain jar
$ws.Range("A7").Value = 'Maven'
$ws.Range("B7").Value = 'Build a jar'
$ws.Range("C7").Value = 'run "mvn package" under the target directory'

# Give the new rows the same "wrap text" formatting used by the existing data rows
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Every data row (2-7) now uses a uniform 36pt row height; the header row keeps its 15pt height
$ws.Range("A2:C7").RowHeight = 36
$ws.Range("A1:C1").RowHeight = 15

# Move the active selection to just past the new last row, like after typing the last entry
[void]$ws.Range("C8").Select()

Write-Output "done"
